# Replace the daily accident report contents of Sheet1 with the next day's
# (2018-11-29) data. The source rows (2-17) get overwritten in place, and five
# additional rows (18-22) are appended, extending the used range to A1:BA22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; A=7; C='No Injuries'; D=35.080462; E=-85.265121; F='2018-11-29'; G='20:51:08'; H='3399 Amnicola Hwy'; K='CHATTANOOGA'; P='20'; Y=11 }
    @{ Row=3; A=12; C='Injuries'; D=35.087713; E=-85.071062; F='2018-11-29'; G='20:29:53'; H='8900 - 8935 Lee Hwy'; K='CHATTANOOGA'; P='20'; Y=11 }
    @{ Row=4; A=17; C='Injuries'; D=35.17114; E=-85.176414; F='2018-11-29'; G='19:09:14'; H='Hixson Pike / Thrasher Pike'; K='HAMILTON COUNTY'; P='19'; Y=11 }
    @{ Row=5; A=43; C='Injuries'; D=35.086343; E=-85.068399; F='2018-11-29'; G='16:47:24'; H='8950 - 8999 Lee Hwy'; K='CHATTANOOGA'; P='16'; Y=11 }
    @{ Row=6; A=44; C='Injuries'; D=35.086343; E=-85.068399; F='2018-11-29'; G='16:47:24'; H='8950 - 8999 Lee Hwy'; K='CHATTANOOGA'; P='16'; Y=11 }
    @{ Row=7; A=54; C='Injuries'; D=35.008065; E=-85.273591; F='2018-11-29'; G='16:14:11'; H='14th Ave / E 29th St'; K='CHATTANOOGA'; P='16'; Y=11 }
    @{ Row=8; A=59; C='Injuries'; D=35.046354; E=-85.278034; F='2018-11-29'; G='16:06:43'; H='600 N Holtzclaw Ave'; K='CHATTANOOGA'; P='16'; Y=11 }
    @{ Row=9; A=63; C='Injuries'; D=35.033649; E=-85.257524; F='2018-11-29'; G='16:01:50'; H='54 SHALLOWFORD RD'; K='CHATTANOOGA'; P='16'; Y=11 }
    @{ Row=10; A=77; C='Injuries'; D=35.145962; E=-85.318082; F='2018-11-29'; G='15:08:56'; H='1607 Anderson Pike'; K='WALDEN'; P='15'; Y=11 }
    @{ Row=11; A=79; C='Unknown Injuries'; D=34.985759; E=-85.226483; F='2018-11-29'; G='15:06:22'; H='Key West Ave / State Line Rd'; K='EAST RIDGE'; P='15'; Y=11 }
    @{ Row=12; A=87; C='Injuries'; D=35.01796; E=-85.142169; F='2018-11-29'; G='14:20:39'; H='1808 Jenkins Rd'; K='CHATTANOOGA'; P='14'; Y=11 }
    @{ Row=13; A=88; C='Injuries'; D=35.01796; E=-85.142169; F='2018-11-29'; G='14:18:48'; H='1808 Jenkins Rd'; K='CHATTANOOGA'; P='14'; Y=11 }
    @{ Row=14; A=97; C='Injuries'; D=35.014497; E=-85.325187; F='2018-11-29'; G='12:17:16'; H='Broad St / Tennessee Ave'; K='CHATTANOOGA'; P='12'; Y=11 }
    @{ Row=15; A=107; C='Injuries'; D=35.024777; E=-85.276914; F='2018-11-29'; G='10:14:25'; H='2106 E Main St'; K='CHATTANOOGA'; P='10'; Y=11 }
    @{ Row=16; A=110; C='Unknown Injuries'; D=35.024527; E=-85.275778; F='2018-11-29'; G='10:02:29'; H='E MAIN ST / S KELLEY ST'; K='CHATTANOOGA'; P='10'; Y=11 }
    @{ Row=17; A=120; C='Injuries'; D=35.182577; E=-85.246697; F='2018-11-29'; G='07:55:04'; H='HIGHWAY 153 / BOY SCOUT RD'; K='CHATTANOOGA'; P='7'; Y=11 }
    @{ Row=18; A=121; C='Injuries'; D=35.182577; E=-85.246697; F='2018-11-29'; G='07:53:17'; H='HIGHWAY 153 / BOY SCOUT RD'; K='CHATTANOOGA'; P='7'; Y=11 }
    @{ Row=19; A=122; C='Unknown Injuries'; D=35.182577; E=-85.246697; F='2018-11-29'; G='07:51:24'; H='Highway 153 / Boy Scout Rd'; K='HAMILTON COUNTY'; P='7'; Y=11 }
    @{ Row=20; A=127; C='Entrapment'; D=35.08104; E=-85.236081; F='2018-11-29'; G='07:13:21'; H='2601 HARRISON PIKE'; K='CHATTANOOGA'; P='7'; Y=11 }
    @{ Row=21; A=129; C='Injuries'; D=35.081593; E=-85.209734; F='2018-11-29'; G='06:52:26'; H='Highway 58 / Bonny Oaks Dr'; K='CHATTANOOGA'; P='6'; Y=11 }
    @{ Row=22; A=130; C='Injuries'; D=35.081593; E=-85.209734; F='2018-11-29'; G='06:51:35'; H='Highway 58 / Bonny Oaks Dr'; K='CHATTANOOGA'; P='6'; Y=11 }
)


foreach ($d in $data) {
    $r = $d.Row

    # New rows (18-22) don't exist yet - create row 18's A cell with the same
    # bold/bordered/centered style used by the rest of column A (copy format
    # from the last existing data row, A17, onto the new A cell only so we
    # don't spray formatting across unused columns).
    if ($r -gt 17) {
        $ws.Cells.Item(17, 1).Copy()
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E

    # Date (F) and Hour (P) look numeric/date-like, so a plain .Value
    # assignment would get auto-coerced into a date serial / number. Route
    # them through a text formula + paste-values so they land as literal
    # shared-string text, matching the source workbook's column typing.
    $ws.Cells.Item($r, 6).Formula = '="' + $d.F + '"'
    $ws.Cells.Item($r, 6).Copy()
    $ws.Cells.Item($r, 6).PasteSpecial(-4163)

    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 11).Value = $d.K

    $ws.Cells.Item($r, 16).Formula = '="' + $d.P + '"'
    $ws.Cells.Item($r, 16).Copy()
    $ws.Cells.Item($r, 16).PasteSpecial(-4163)

    $ws.Cells.Item($r, 25).Value = $d.Y
}

$excel.CutCopyMode = 0
